$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a rolling "MarketBeat rank" watch report: column B always
# holds the most-recent scan date, with older scans pushed right each time a
# new scan is appended. This run ("10th") adds three new date columns
# (two for Jun_26 and one for Jun_27) in front of the existing data, and
# appends a brand-new analyst group (two new rows) at the bottom.
# ---------------------------------------------------------------------------

# 1) Make room for the new scan columns: insert 3 blank columns starting at B
#    (old B,C,D,E - the four scan-date columns - shift right to E,F,G,H).
$ws.Columns("B:D").Insert()

# Re-apply the custom 8-char width to every date column (C..H) - Insert()
# brings width=8 along but drops the "custom width" flag.
$ws.Columns("C:H").ColumnWidth = 7.1666667

# 2) New header row values for the freshly inserted columns; E1:H1 already
#    hold the correct (shifted) former B1:E1 values.
$ws.Range("B1").Value = "Jun_27"
$ws.Range("C1").Value = "Jun_26"
$ws.Range("D1").Value = "Jun_26"

# 3) Default every new cell in B:D (rows 2-27) to "UN", matching the rest of
#    the sheet; row-specific alerts are overwritten afterwards.
for ($r = 2; $r -le 27; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# 4) Row 20 (Bank of America): new alert on 6/25/2018 lands in the Jun_26
#    columns (C20 & D20).
$ws.Range("C20").Value = '6/25/2018,Initiates,Buy,$420.00'
$ws.Range("D20").Value = '6/25/2018,Initiates,Buy,$420.00'

# 5) Row 27 (Barclays): new alert on 6/21/2018 lands in the Jun_26 columns
#    (C27 & D27); C27 is highlighted the same way the sheet already
#    highlights the other "fresh alert" cells (E6/E7, now H6/H7).
$ws.Range("C27").Value = '6/21/2018,Raises Target,Overweight,$410.00'
$ws.Range("D27").Value = '6/21/2018,Raises Target,Overweight,$410.00'
$ws.Range("C27").Interior.Color = $ws.Range("H6").Interior.Color

# 6) Append the new analyst group at the bottom of the table.
$ws.Range("A28").Value = "Benchmark"
$ws.Range("B28").Value = "UN"
$ws.Range("C28").Value = "UN"
$ws.Range("D28").Value = "UN"

$ws.Range("A29").Value = "Evercore ISI"
$ws.Range("B29").Value = "UN"
$ws.Range("C29").Value = "UN"
$ws.Range("D29").Value = "UN"
